$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Monday hours for the week commencing 43178 (row 10) from 4.25 to 8.25
$ws.Range("B10").Value = 8.25

# Update the active cell selection to B11 as in the edited file
$ws.Range("B11").Select()
